# Weekly fruit/vegetable price update: add one new week's worth of records
# (two rows: "Copenhague"/"Primera" and "Crespo record"/"Primera") at the
# top of the data block, pushing all existing data rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the first data row (row 611),
# shifting every existing record down by two rows.
$ws.Rows("611:612").Insert()

# --- New row 611: Copenhague / Primera ---
$ws.Range("A611").Value = 4
$ws.Range("B611").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C611").Value = "Los Lagos"
$ws.Range("D611").Value2 = 45034
$ws.Range("D611").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E611").Value = 10
$ws.Range("F611").Value = 100112006
$ws.Range("G611").Value = "Repollo"
$ws.Range("H611").Value = "Copenhague"
$ws.Range("I611").Value = "Primera"
$ws.Range("J611").Value = 700
$ws.Range("K611").Value = 2000
$ws.Range("L611").Value = 2000
$ws.Range("M611").Value = 2000
$ws.Range("N611").Value = "$/unidad"
$ws.Range("O611").Value = "Región Metropolitana"
$ws.Range("P611").Value = 2000
$ws.Range("Q611").Value = 1
$ws.Range("R611").Value = "Hortaliza"

# --- New row 612: Crespo record / Primera ---
$ws.Range("A612").Value = 4
$ws.Range("B612").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C612").Value = "Los Lagos"
$ws.Range("D612").Value2 = 45034
$ws.Range("D612").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E612").Value = 10
$ws.Range("F612").Value = 100112006
$ws.Range("G612").Value = "Repollo"
$ws.Range("H612").Value = "Crespo record"
$ws.Range("I612").Value = "Primera"
$ws.Range("J612").Value = 700
$ws.Range("K612").Value = 2000
$ws.Range("L612").Value = 2000
$ws.Range("M612").Value = 2000
$ws.Range("N612").Value = "$/unidad"
$ws.Range("O612").Value = "Región Metropolitana"
$ws.Range("P612").Value = 2000
$ws.Range("Q612").Value = 1
$ws.Range("R612").Value = "Hortaliza"
